$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.605581045150757
$ws.Range("B1").Value = 2.503911733627319
$ws.Range("C1").Value = 2.017632961273193
$ws.Range("D1").Value = 2.018800020217896
$ws.Range("E1").Value = 2.256995439529419
